# issue #5: property land done
# Rewrites the "土地" (land) sheet into the flattened export schema
# (name/area/share_portion/... + property_category/category/date/
# legislator_name/legislator_id/source_file/index) and normalises stray
# whitespace / punctuation in the "建物" (building) and "汽車" (car) sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "土地" (land) — new header row + extra trailing columns I:O
# ---------------------------------------------------------------------
$land = $wb.Worksheets.Item("土地")

# Make the new columns I:O inherit the existing header / data-row look
# (same border+bold style as B1:H1, same plain style as B2:H2) before
# writing values into them.
$land.Range("H1").Copy()
$land.Range("I1:O1").PasteSpecial(-4122)
$land.Range("H2").Copy()
$land.Range("I2:O4").PasteSpecial(-4122)

# K2:K4 hold an ISO date string ("2012-05-01"); force text formatting
# first so Excel doesn't silently reinterpret it as a date serial.
$land.Range("K2:K4").NumberFormat = "@"

# Header row
$land.Range("B1").Value = "name"
$land.Range("C1").Value = "area"
$land.Range("D1").Value = "share_portion"
$land.Range("E1").Value = "owner"
$land.Range("F1").Value = "register_date"
$land.Range("G1").Value = "register_reason"
$land.Range("H1").Value = "acquire_value"
$land.Range("I1").Value = "property_category"
$land.Range("J1").Value = "category"
$land.Range("K1").Value = "date"
$land.Range("L1").Value = "legislator_name"
$land.Range("M1").Value = "legislator_id"
$land.Range("N1").Value = "source_file"
$land.Range("O1").Value = "index"

# Row 2
$land.Range("B2").Value = "臺北市大安區仁愛段二小段06010000地號"
$land.Range("C2").Value = 509
$land.Range("D2").Value = "30000分之1703"
$land.Range("E2").Value = "林世嘉"
$land.Range("F2").Value = "93年01月02日"
$land.Range("G2").Value = "買賣"
$land.Range("H2").Value = "(超過五年）"
$land.Range("I2").Value = "land"
$land.Range("J2").Value = "normal"
$land.Range("K2").Value = "2012-05-01"
$land.Range("L2").Value = "林世嘉"
$land.Range("M2").Value = 1740
$land.Range("N2").Value = "tmpada11"
$land.Range("O2").Value = 14

# Row 3
$land.Range("B3").Value = "臺北市中山區長安段四小段02980000地號"
$land.Range("C3").Value = 1069
$land.Range("D3").Value = "10000分之211"
$land.Range("E3").Value = "蔡篤堅"
$land.Range("F3").Value = "92年01月17日"
$land.Range("G3").Value = "買賣"
$land.Range("H3").Value = "(超過五年）"
$land.Range("I3").Value = "land"
$land.Range("J3").Value = "normal"
$land.Range("K3").Value = "2012-05-01"
$land.Range("L3").Value = "林世嘉"
$land.Range("M3").Value = 1740
$land.Range("N3").Value = "tmpada11"
$land.Range("O3").Value = 15

# Row 4
$land.Range("B4").Value = "臺北市中山區長安段四小段02300000地號"
$land.Range("C4").Value = 448
$land.Range("D4").Value = "10000分之180"
$land.Range("E4").Value = "蔡篤堅"
$land.Range("F4").Value = "97年04月03日"
$land.Range("G4").Value = "買賣"
$land.Range("H4").Value = "15000000(房地總價額）"
$land.Range("I4").Value = "land"
$land.Range("J4").Value = "normal"
$land.Range("K4").Value = "2012-05-01"
$land.Range("L4").Value = "林世嘉"
$land.Range("M4").Value = 1740
$land.Range("N4").Value = "tmpada11"
$land.Range("O4").Value = 16

# Drop the now-unneeded text format override so K2:K4 fall back to the
# sheet's plain/default look (matches the rest of the row).
$land.Range("K2:K4").ClearFormats()

# ---------------------------------------------------------------------
# Sheet "建物" (building) — strip stray spaces / hyphens from values
# ---------------------------------------------------------------------
$building = $wb.Worksheets.Item("建物")

$building.Range("B2").Value = "臺北市大安區仁愛段二小段05345000建號"
$building.Range("F2").Value = "93年01月02日"

$building.Range("B3").Value = "臺北市大安區仁愛段二小段03854000建號"
$building.Range("F3").Value = "93年01月02日"

$building.Range("B4").Value = "臺北市中山區長安段四小段01678000建號"
$building.Range("F4").Value = "92年01月17曰"

$building.Range("B5").Value = "臺北市中山區長安段四小段02074000建號"
$building.Range("D5").Value = "10000分之197"
$building.Range("F5").Value = "97年04月03日"
$building.Range("H5").Value = "15000000(同上筆土地）"

$building.Range("B6").Value = "臺北市中山區長安段四小段03734000建號"
$building.Range("F6").Value = "97年04月03日"
$building.Range("H6").Value = "15000000(同上筆土地）"

# ---------------------------------------------------------------------
# Sheet "汽車" (car) — strip stray space from the model name
# ---------------------------------------------------------------------
$car = $wb.Worksheets.Item("汽車")
$car.Range("B2").Value = "VOLKSWAGENPSDDAT2.0"
